$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 209
$ws1.Range("F4").Value = 12845
$ws1.Range("F10").Value = 216
$ws1.Range("F11").Value = 465
$ws1.Range("F23").Value = 117

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 209
$ws4.Range("F4").Value = 12845
$ws4.Range("F10").Value = 216
$ws4.Range("F11").Value = 465
$ws4.Range("F23").Value = 117
